$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 161.83333
$ws.Range("I5").Value = 161.83333
$ws.Range("K5").Value = 161.83333
$ws.Range("M5").Value = -46.83332999999999
$ws.Range("H33").Value = 399.875
$ws.Range("I33").Value = 399.875
$ws.Range("K33").Value = 399.875
$ws.Range("M33").Value = -170.875
$ws.Range("H55").Value = 499.36365
$ws.Range("J55").Value = 492
$ws.Range("L55").Value = 492
$ws.Range("N55").Value = -920
$ws.Range("H74").Value = 11827.52
$ws.Range("I74").Value = 12051.652
$ws.Range("K74").Value = 12051.652
$ws.Range("M74").Value = -11115.652
$ws.Range("H76").Value = 3999.5
$ws.Range("I76").Value = 3999.3333
$ws.Range("K76").Value = 3999.3333
$ws.Range("M76").Value = -3684.3333
$ws.Range("H77").Value = 11827.52
$ws.Range("I77").Value = 12051.652
$ws.Range("K77").Value = 60258.26
$ws.Range("M77").Value = -55578.26
$ws.Range("H79").Value = 3999.5
$ws.Range("I79").Value = 3999.3333
$ws.Range("K79").Value = 3999.3333
$ws.Range("M79").Value = -2907.3333
$ws.Range("H107").Value = 1220.1333
$ws.Range("J107").Value = 1679.6
$ws.Range("L107").Value = 1679.6
$ws.Range("N107").Value = -5519.6
$ws.Range("H111").Value = 4505.4443
$ws.Range("I111").Value = 6091.6665
$ws.Range("J111").Value = 1333
$ws.Range("K111").Value = 18274.9995
$ws.Range("L111").Value = 3999
$ws.Range("M111").Value = -15207.9995
$ws.Range("N111").Value = -10133
$ws.Range("H112").Value = 1918.4445
$ws.Range("J112").Value = 2108.3125
$ws.Range("L112").Value = 6324.9375
$ws.Range("N112").Value = -8540.9375
$ws.Range("H121").Value = 3000
$ws.Range("J121").Value = 3000
$ws.Range("L121").Value = 9000
$ws.Range("N121").Value = -12494
$ws.Range("H127").Value = 650.75
$ws.Range("I127").Value = 650.75
$ws.Range("K127").Value = 1952.25
$ws.Range("M127").Value = 3007.75
$ws.Range("H137").Value = 8335867
$ws.Range("I137").Value = 8335867
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 25007601
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -25005051
$ws.Range("N137").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2366
$ws.Range("I2").Value = 2308.3333
$ws.Range("J2").Value = 2464.8572
$ws.Range("K2").Value = 2308.3333
$ws.Range("L2").Value = 2464.8572
$ws.Range("M2").Value = -2195.3333
$ws.Range("N2").Value = -2690.8572
$ws.Range("H32").Value = 1883110.6
$ws.Range("I32").Value = 843657.5600000001
$ws.Range("K32").Value = 843657.5600000001
$ws.Range("M32").Value = -843370.5600000001
$ws.Range("H55").Value = 63495.5
$ws.Range("J55").Value = 76976.5
$ws.Range("L55").Value = 76976.5
$ws.Range("N55").Value = -77606.5
$ws.Range("H63").Value = 5214.4287
$ws.Range("I63").Value = 2375.75
$ws.Range("J63").Value = 8999.333000000001
$ws.Range("K63").Value = 2375.75
$ws.Range("L63").Value = 8999.333000000001
$ws.Range("M63").Value = -1689.75
$ws.Range("N63").Value = -10371.333
$ws.Range("H66").Value = 5214.4287
$ws.Range("I66").Value = 2375.75
$ws.Range("J66").Value = 8999.333000000001
$ws.Range("K66").Value = 11878.75
$ws.Range("L66").Value = 44996.665
$ws.Range("M66").Value = -8446.75
$ws.Range("N66").Value = -51860.665
$ws.Range("H74").Value = 50279904
$ws.Range("I74").Value = 506594.9
$ws.Range("J74").Value = 111113944
$ws.Range("K74").Value = 506594.9
$ws.Range("L74").Value = 111113944
$ws.Range("M74").Value = -505720.9
$ws.Range("N74").Value = -111115692
$ws.Range("H77").Value = 50279904
$ws.Range("I77").Value = 506594.9
$ws.Range("J77").Value = 111113944
$ws.Range("K77").Value = 2532974.5
$ws.Range("L77").Value = 555569720
$ws.Range("M77").Value = -2528606.5
$ws.Range("N77").Value = -555578456
$ws.Range("H116").Value = 2366
$ws.Range("I116").Value = 2308.3333
$ws.Range("J116").Value = 2464.8572
$ws.Range("K116").Value = 2308.3333
$ws.Range("L116").Value = 2464.8572
$ws.Range("M116").Value = -14.33329999999978
$ws.Range("N116").Value = -7052.8572

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2366
$ws.Range("I3").Value = 2308.3333
$ws.Range("J3").Value = 2464.8572
$ws.Range("K3").Value = 2308.3333
$ws.Range("L3").Value = 2464.8572
$ws.Range("M3").Value = -2194.3333
$ws.Range("N3").Value = -2692.8572
$ws.Range("H35").Value = 57074
$ws.Range("J35").Value = 57074
$ws.Range("L35").Value = 57074
$ws.Range("N35").Value = -57694
$ws.Range("H86").Value = 2716.8
$ws.Range("I86").Value = 2678.4348
$ws.Range("J86").Value = 2842.8572
$ws.Range("K86").Value = 2678.4348
$ws.Range("L86").Value = 2842.8572
$ws.Range("M86").Value = -1555.4348
$ws.Range("N86").Value = -5088.8572
$ws.Range("H89").Value = 2716.8
$ws.Range("I89").Value = 2678.4348
$ws.Range("J89").Value = 2842.8572
$ws.Range("K89").Value = 13392.174
$ws.Range("L89").Value = 14214.286
$ws.Range("M89").Value = -7776.173999999999
$ws.Range("N89").Value = -25446.286
$ws.Range("H134").Value = 4030.4614
$ws.Range("I134").Value = 4488.4443
$ws.Range("K134").Value = 13465.3329
$ws.Range("M134").Value = -10930.3329

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2913042
$ws.Range("I31").Value = 1898.8125
$ws.Range("J31").Value = 4638163.5
$ws.Range("K31").Value = 1898.8125
$ws.Range("L31").Value = 4638163.5
$ws.Range("M31").Value = -1603.8125
$ws.Range("N31").Value = -4638753.5
$ws.Range("H34").Value = 2913042
$ws.Range("I34").Value = 1898.8125
$ws.Range("J34").Value = 4638163.5
$ws.Range("K34").Value = 1898.8125
$ws.Range("L34").Value = 4638163.5
$ws.Range("M34").Value = -1696.8125
$ws.Range("N34").Value = -4638567.5
$ws.Range("H97").Value = 41597
$ws.Range("J97").Value = 41597
$ws.Range("L97").Value = 41597
$ws.Range("N97").Value = -43579
$ws.Range("H132").Value = 4378.0527
$ws.Range("I132").Value = 3093.923
$ws.Range("K132").Value = 9281.769
$ws.Range("M132").Value = -6751.769
$ws.Range("H134").Value = 3753.2
$ws.Range("I134").Value = 3959.4
$ws.Range("K134").Value = 11878.2
$ws.Range("M134").Value = -9343.200000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 43.8
$ws.Range("J12").Value = 43.8
$ws.Range("L12").Value = 131.4
$ws.Range("N12").Value = -477.4
$ws.Range("H34").Value = 439.66666
$ws.Range("I34").Value = 313.33334
$ws.Range("K34").Value = 940.0000200000001
$ws.Range("M34").Value = -856.0000200000001
$ws.Range("H113").Value = 1068
$ws.Range("I113").Value = 876.75
$ws.Range("K113").Value = 2630.25
$ws.Range("M113").Value = -460.25
$ws.Range("H129").Value = 168.14285
$ws.Range("I129").Value = 168.14285
$ws.Range("K129").Value = 504.42855
$ws.Range("M129").Value = 4495.57145

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1721.8286
$ws.Range("I132").Value = 1638.303
$ws.Range("J132").Value = 3100
$ws.Range("K132").Value = 4914.909000000001
$ws.Range("L132").Value = 9300
$ws.Range("M132").Value = -2384.909000000001
$ws.Range("N132").Value = -14360

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 57305.78
$ws.Range("I40").Value = 62096.81
$ws.Range("K40").Value = 62096.81
$ws.Range("M40").Value = -61960.81
$ws.Range("H93").Value = 1700.091
$ws.Range("I93").Value = 1528.56
$ws.Range("J93").Value = 2236.125
$ws.Range("K93").Value = 1528.56
$ws.Range("L93").Value = 2236.125
$ws.Range("M93").Value = -280.5599999999999
$ws.Range("N93").Value = -4732.125
$ws.Range("H122").Value = 3479.2
$ws.Range("I122").Value = 3479.2
$ws.Range("K122").Value = 10437.6
$ws.Range("M122").Value = -7987.599999999999
$ws.Range("H132").Value = 6664.2
$ws.Range("I132").Value = 3495
$ws.Range("K132").Value = 10485
$ws.Range("M132").Value = -7955

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 15273
$ws.Range("I52").Value = 16030.667
$ws.Range("K52").Value = 16030.667
$ws.Range("M52").Value = -15804.667
$ws.Range("H107").Value = 482.41177
$ws.Range("I107").Value = 388.6
$ws.Range("K107").Value = 1165.8
$ws.Range("M107").Value = 754.1999999999998
$ws.Range("H126").Value = 8424.933999999999
$ws.Range("I126").Value = 8424.933999999999
$ws.Range("K126").Value = 25274.802
$ws.Range("M126").Value = -22804.802
